# Generate Report for Archive
#
# The localization-status report is regenerated and the row belonging to
# 5b362aeb-48e6-485c-b12c-2f29b618fa7a moves up from the 5th data row to the
# 3rd data row (worksheet row 6 -> row 4), pushing the 5e78426a and
# 39c22f56 rows down by one each. This happens identically on all three
# worksheets (Overview, zh-cn, de-de); the zh-cn/de-de sheets also carry the
# "Latest Handoff File" / "Latest Handoff Datetime" columns for that row.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Cells.Item(4,1).Value = "5b362aeb-48e6-485c-b12c-2f29b618fa7a.md"
$ws.Cells.Item(4,2).Value = "Ready for handoff"
$ws.Cells.Item(4,3).Value = "Ready for handoff"

$ws.Cells.Item(5,1).Value = "5e78426a-9986-44d8-845a-7bbcd9961e21.md"
$ws.Cells.Item(5,2).Value = "In Translation"
$ws.Cells.Item(5,3).Value = "In Translation"

$ws.Cells.Item(6,1).Value = "39c22f56-0b82-4ddb-ac83-a2c10f6fd190.md"
$ws.Cells.Item(6,2).Value = "Ready for handoff"
$ws.Cells.Item(6,3).Value = "Ready for handoff"

foreach ($hl in $ws.Hyperlinks) {
    $r = $hl.Range.Row
    if ($r -eq 4) {
        $hl.TextToDisplay = "5b362aeb-48e6-485c-b12c-2f29b618fa7a.md"
    } elseif ($r -eq 5) {
        $hl.TextToDisplay = "5e78426a-9986-44d8-845a-7bbcd9961e21.md"
    } elseif ($r -eq 6) {
        $hl.TextToDisplay = "39c22f56-0b82-4ddb-ac83-a2c10f6fd190.md"
    }
}

# ---- zh-cn sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Cells.Item(4,1).Value = "5b362aeb-48e6-485c-b12c-2f29b618fa7a.md"
$ws.Cells.Item(4,2).Value = "Ready for handoff"
$ws.Cells.Item(4,3).Value = "5b362aeb-48e6-485c-b12c-2f29b618fa7a.0dcaf32b3a5cd8d143d8faa192f351c4872822cb.zh-cn.xlf"
$ws.Cells.Item(4,4).Value = "2016-02-29 13:06:35"

$ws.Cells.Item(5,1).Value = "5e78426a-9986-44d8-845a-7bbcd9961e21.md"
$ws.Cells.Item(5,2).Value = "In Translation"
$ws.Cells.Item(5,3).Value = "5e78426a-9986-44d8-845a-7bbcd9961e21.e30fb42c2d185b3f269b21f9cdc9c0498a59b7a4.zh-cn.xlf"
$ws.Cells.Item(5,4).Value = "2016-02-29 12:54:22"

$ws.Cells.Item(6,1).Value = "39c22f56-0b82-4ddb-ac83-a2c10f6fd190.md"
$ws.Cells.Item(6,2).Value = "Ready for handoff"
$ws.Cells.Item(6,3).Value = "39c22f56-0b82-4ddb-ac83-a2c10f6fd190.9b93316d6c03bfaf63c37925043c36f5b4b5424a.zh-cn.xlf"
$ws.Cells.Item(6,4).Value = "2016-02-29 12:57:17"

foreach ($hl in $ws.Hyperlinks) {
    $r = $hl.Range.Row
    $c = $hl.Range.Column
    if ($r -eq 4 -and $c -eq 1) {
        $hl.TextToDisplay = "5b362aeb-48e6-485c-b12c-2f29b618fa7a.md"
    } elseif ($r -eq 4 -and $c -eq 3) {
        $hl.TextToDisplay = "5b362aeb-48e6-485c-b12c-2f29b618fa7a.0dcaf32b3a5cd8d143d8faa192f351c4872822cb.zh-cn.xlf"
    } elseif ($r -eq 5 -and $c -eq 1) {
        $hl.TextToDisplay = "5e78426a-9986-44d8-845a-7bbcd9961e21.md"
    } elseif ($r -eq 5 -and $c -eq 3) {
        $hl.TextToDisplay = "5e78426a-9986-44d8-845a-7bbcd9961e21.e30fb42c2d185b3f269b21f9cdc9c0498a59b7a4.zh-cn.xlf"
    } elseif ($r -eq 6 -and $c -eq 1) {
        $hl.TextToDisplay = "39c22f56-0b82-4ddb-ac83-a2c10f6fd190.md"
    } elseif ($r -eq 6 -and $c -eq 3) {
        $hl.TextToDisplay = "39c22f56-0b82-4ddb-ac83-a2c10f6fd190.9b93316d6c03bfaf63c37925043c36f5b4b5424a.zh-cn.xlf"
    }
}

# ---- de-de sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Cells.Item(4,1).Value = "5b362aeb-48e6-485c-b12c-2f29b618fa7a.md"
$ws.Cells.Item(4,2).Value = "Ready for handoff"
$ws.Cells.Item(4,3).Value = "5b362aeb-48e6-485c-b12c-2f29b618fa7a.0dcaf32b3a5cd8d143d8faa192f351c4872822cb.de-de.xlf"
$ws.Cells.Item(4,4).Value = "2016-02-29 13:06:48"

$ws.Cells.Item(5,1).Value = "5e78426a-9986-44d8-845a-7bbcd9961e21.md"
$ws.Cells.Item(5,2).Value = "In Translation"
$ws.Cells.Item(5,3).Value = "5e78426a-9986-44d8-845a-7bbcd9961e21.e30fb42c2d185b3f269b21f9cdc9c0498a59b7a4.de-de.xlf"
$ws.Cells.Item(5,4).Value = "2016-02-29 12:54:33"

$ws.Cells.Item(6,1).Value = "39c22f56-0b82-4ddb-ac83-a2c10f6fd190.md"
$ws.Cells.Item(6,2).Value = "Ready for handoff"
$ws.Cells.Item(6,3).Value = "39c22f56-0b82-4ddb-ac83-a2c10f6fd190.9b93316d6c03bfaf63c37925043c36f5b4b5424a.de-de.xlf"
$ws.Cells.Item(6,4).Value = "2016-02-29 12:57:30"

foreach ($hl in $ws.Hyperlinks) {
    $r = $hl.Range.Row
    $c = $hl.Range.Column
    if ($r -eq 4 -and $c -eq 1) {
        $hl.TextToDisplay = "5b362aeb-48e6-485c-b12c-2f29b618fa7a.md"
    } elseif ($r -eq 4 -and $c -eq 3) {
        $hl.TextToDisplay = "5b362aeb-48e6-485c-b12c-2f29b618fa7a.0dcaf32b3a5cd8d143d8faa192f351c4872822cb.de-de.xlf"
    } elseif ($r -eq 5 -and $c -eq 1) {
        $hl.TextToDisplay = "5e78426a-9986-44d8-845a-7bbcd9961e21.md"
    } elseif ($r -eq 5 -and $c -eq 3) {
        $hl.TextToDisplay = "5e78426a-9986-44d8-845a-7bbcd9961e21.e30fb42c2d185b3f269b21f9cdc9c0498a59b7a4.de-de.xlf"
    } elseif ($r -eq 6 -and $c -eq 1) {
        $hl.TextToDisplay = "39c22f56-0b82-4ddb-ac83-a2c10f6fd190.md"
    } elseif ($r -eq 6 -and $c -eq 3) {
        $hl.TextToDisplay = "39c22f56-0b82-4ddb-ac83-a2c10f6fd190.9b93316d6c03bfaf63c37925043c36f5b4b5424a.de-de.xlf"
    }
}
